# BOM update: C4 merged into C1 (same 0.1uF cap), Single layer pad TP row
# removed (pcb change), and R4 resistor value changed from 49.9K to 200K
# (manufacturer part number cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("CAP CER 0.068UF 16V X7R 0603" / C4 / CC0603KRX7R7BB683) is
# removed - C4 now shares the C1 part (0.1uF cap), so bump C1's row to
# cover both designators and quantity 2, then delete the old C4 row.
$ws.Range("B2").Value = "C1, C4"
$ws.Range("C2").Value = 2
$ws.Rows(4).Delete()

# "Single layer pad TP" / "JP2 Vref, JP3 Vout" row (now row 6 after the
# shift above) is removed entirely - pcb no longer has those test points.
$ws.Rows(6).Delete()

# R4 changed from RES SMD 49.9K OHM 1% 1/10W 0603 to RES SMD 200K OHM 1%
# 1/10W 0603, and its old manufacturer part number (RC0603FR-0749K9L) is
# cleared since the new part has none listed. After the two deletions
# above, the R4 row is now row 10.
$ws.Range("A10").Value = "RES SMD 200K OHM 1% 1/10W 0603"
$ws.Range("D10").ClearContents()
